$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.117.35"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "3.107.87"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.17"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.76"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.49"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.153"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.476"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000247"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.67"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "3.622.34"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").Value = "67.082.97"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.08"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").Value = "3.104.34"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "490.84"
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.84"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.700"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "83.78"
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.04"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("E25").Value = "  -3.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.53"
$ws.Range("E26").Value = "  +3.97%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.88"
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.22"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").Value = "0.0₃0933"
$ws.Range("E33").Value = "  -7.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.972"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "46.83"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("E41").Value = "  -2.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "385.00"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "2.797.31"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  -9.01%  "
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.02"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.92"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.69"
$ws.Range("E51").Value = "  -2.06%  "
